$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select() | Out-Null
